# Fruta / hortaliza, semanal
# Insert two new weekly rows (current row 33/34, shifting the old
# rows 33-49 down to 35-51) for the Arándano (blue) - Terminal
# Hortofrutícola Agro Chillán sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 33 (pushes old row 33
# down to row 35, etc.), mirroring how Excel's Rows.Insert behaves -
# the inserted rows inherit formatting (e.g. the date style on column D)
# from the surrounding rows.
$ws.Rows.Item(33).Insert()
$ws.Rows.Item(33).Insert()

# New row 33: Primera quality, week of 2023-12-04 (serial 45264)
$ws.Cells.Item(33, 1).Value = 7
$ws.Cells.Item(33, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(33, 3).Value = "Ñuble"
$ws.Cells.Item(33, 4).Value = 45264
$ws.Cells.Item(33, 5).Value = 16
$ws.Cells.Item(33, 6).Value = "Fruta"
$ws.Cells.Item(33, 7).Value = 100101
$ws.Cells.Item(33, 8).Value = "Berries"
$ws.Cells.Item(33, 9).Value = 100101001
$ws.Cells.Item(33, 10).Value = "Arándano (blue)"
$ws.Cells.Item(33, 11).Value = "Sin especificar"
$ws.Cells.Item(33, 12).Value = "Primera"
$ws.Cells.Item(33, 13).Value = 60
$ws.Cells.Item(33, 14).Value = 6000
$ws.Cells.Item(33, 15).Value = 6000
$ws.Cells.Item(33, 16).Value = 6000
$ws.Cells.Item(33, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(33, 18).Value = "Región de Ñuble"
$ws.Cells.Item(33, 19).Value = 3000
$ws.Cells.Item(33, 20).Value = 2

# New row 34: Segunda quality, same week (serial 45264)
$ws.Cells.Item(34, 1).Value = 7
$ws.Cells.Item(34, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(34, 3).Value = "Ñuble"
$ws.Cells.Item(34, 4).Value = 45264
$ws.Cells.Item(34, 5).Value = 16
$ws.Cells.Item(34, 6).Value = "Fruta"
$ws.Cells.Item(34, 7).Value = 100101
$ws.Cells.Item(34, 8).Value = "Berries"
$ws.Cells.Item(34, 9).Value = 100101001
$ws.Cells.Item(34, 10).Value = "Arándano (blue)"
$ws.Cells.Item(34, 11).Value = "Sin especificar"
$ws.Cells.Item(34, 12).Value = "Segunda"
$ws.Cells.Item(34, 13).Value = 60
$ws.Cells.Item(34, 14).Value = 5000
$ws.Cells.Item(34, 15).Value = 5000
$ws.Cells.Item(34, 16).Value = 5000
$ws.Cells.Item(34, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(34, 18).Value = "Región de Ñuble"
$ws.Cells.Item(34, 19).Value = 2500
$ws.Cells.Item(34, 20).Value = 2

# Make sure column D keeps its date number-format style for the two new rows
$ws.Cells.Item(33, 4).Style = $ws.Cells.Item(35, 4).Style
$ws.Cells.Item(34, 4).Style = $ws.Cells.Item(35, 4).Style
